$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign a value that "looks like" a number but must be stored as text
# (matches the source data, which keeps significance-star annotations like
# "0.82*" as text). Briefly mark the cell as Text-formatted so Excel does
# not auto-convert the literal to a number, then restore the default
# "Normal" style so the cell's style/format stays as it was originally.
function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "B2" "0.17"
Set-TextValue "B3" "-0.01"
Set-TextValue "B4" "-0.09"

Set-TextValue "C2" "44.29***"
Set-TextValue "C3" "2.21***"
Set-TextValue "C4" "0.98"

Set-TextValue "D2" "-0.89"
Set-TextValue "D3" "0.46***"
Set-TextValue "D4" "0.82*"
